# Apply content + formatting updates to match the target revision of sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Cell value updates ---
$ws.Range("B10").Value = '198273 - Domingos Savio Giordani'
$ws.Range("C10").Value = '198273 - Domingos Savio Giordani'
$ws.Range("A13").Value = 'Programa resumido:'
$ws.Range("B13").Value = 'Semestral'
$ws.Range("C13").Value = 'Semestral'
$ws.Range("A14").Value = 'Short syllabus:'
$ws.Range("B14").Value = 'Training and team work, Communication, Systematic Innovation, Legislation, Project Management. Problem Identification, Project Formulation, Problem Specification, Available Knowledge Analysis, Evaluation and Decision Making, Timeline, Reporting, Project Submission'
$ws.Range("C14").Value = 'Training and team work, Communication, Systematic Innovation, Legislation, Project Management. Problem Identification, Project Formulation, Problem Specification, Available Knowledge Analysis, Evaluation and Decision Making, Timeline, Reporting, Project Submission'
$ws.Range("A15").Value = 'Programa:'
$ws.Range("A16").Value = 'Syllabus:'
$ws.Range("B16").Value = 'Training and work in teams and communication - the development of skills essential to work in teams; Systematic Innovation - development of innovative solutions, systematization and characteristics; Legislation - notions of legislation applied to corporate action; Project Management and Schedule - Methodologies and necessary schematizations with the management elements; Problem Identification - systematization of actions to locate causes; Formulation of the Project - presentation of the managerial aspects necessary for the development of the project, Management Plan, Project Analytical Structure (EAP) etc; Specification of Problems - systematization of problems within the areas of knowledge; Analysis of Available Knowledge, Evaluation and Decision Making; Reporting - formatting within ABNT standards; Presentation of Projects.'
$ws.Range("C16").Value = 'Training and work in teams and communication - the development of skills essential to work in teams; Systematic Innovation - development of innovative solutions, systematization and characteristics; Legislation - notions of legislation applied to corporate action; Project Management and Schedule - Methodologies and necessary schematizations with the management elements; Problem Identification - systematization of actions to locate causes; Formulation of the Project - presentation of the managerial aspects necessary for the development of the project, Management Plan, Project Analytical Structure (EAP) etc; Specification of Problems - systematization of problems within the areas of knowledge; Analysis of Available Knowledge, Evaluation and Decision Making; Reporting - formatting within ABNT standards; Presentation of Projects.'
$ws.Range("A17").Value = 'Avaliação:'
$ws.Range("A18").Value = 'Método:'

# B18 is a brand-new cell (row 18 previously only had column A). Column B is
# covered by two overlapping <col> style rules (style 1, then overridden to
# style 2); pull the correct sibling format across explicitly before writing
# the value so B18 ends up styled like the rest of column B (s="2") instead
# of inheriting the first, less-specific column rule.
$ws.Range("B19").Copy()
$ws.Range("B18").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$ws.Range("B18").Value = '198273 - Domingos Savio Giordani'
$ws.Range("C18").Value = '198273 - Domingos Savio Giordani'
$ws.Range("A19").Value = 'Critério:'
$ws.Range("A20").Value = 'Norma de recuperação:'
$ws.Range("A21").Value = 'Bibliografia:'

# B15/C15 become the text "01/01/2020" (must stay text, not be parsed as a date,
# so copy the already-text value from B8 instead of assigning a literal string).
$ws.Range("B8").Copy()
$ws.Range("B15").PasteSpecial(-4163)
$ws.Range("B8").Copy()
$ws.Range("C15").PasteSpecial(-4163)
$excel.CutCopyMode = 0

# --- Cell clears ---
$ws.Range("B17").Clear()
$ws.Range("C17").Clear()

# --- Row height updates ---
$ws.Rows.Item(13).RowHeight = 60
$ws.Rows.Item(15).RowHeight = 120
$ws.Rows.Item(17).AutoFit()
$ws.Rows.Item(18).RowHeight = 60
$ws.Rows.Item(21).RowHeight = 120

# --- Remove the now-deleted last row (old row 22, Bibliografia block) ---
$ws.Rows.Item(22).Delete()
